$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (TimeStamp, CC_Régularisations, DEB_Récurrent, DEB_Trans,
# ENC_Détails, ENC_Entête, FAC_Comptes_Clients, FAC_Détails, FAC_Entête,
# FAC_Projets_Détails, FAC_Projets_Entête, FAC_Sommaire_Taux, GL_EJ_Récurrente,
# GL_Trans, TEC_Local)
$newRows = @(
    @(45700.735578703701, 8, 6, 211, 386, 362, 388, 2681, 388, 1216, 119, 304, 30, 3087, 4051),
    @(45700.739432870374, 8, 6, 211, 386, 362, 388, 2681, 388, 1216, 119, 304, 30, 3087, 4051),
    @(45700.743564814817, 8, 6, 211, 386, 362, 388, 2681, 388, 1216, 119, 304, 30, 3087, 4051),
    @(45700.921655092592, 8, 6, 215, 389, 365, 388, 2681, 388, 1216, 119, 304, 30, 3105, 4077),
    @(45700.931967592594, 8, 6, 215, 389, 365, 388, 2681, 388, 1216, 119, 304, 30, 3105, 4077),
    @(45700.935011574074, 8, 6, 215, 389, 365, 388, 2681, 388, 1216, 119, 304, 30, 3105, 4077)
)

$startRow = 47
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $newRows[$i]
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowValues[$col - 1]
    }
}
